# Add the new "productimages" worksheet after the existing "register" sheet,
# populate it with the product/image-count test data, and format it to match
# the "register" sheet's look (bordered, wrapped cells; bold yellow header).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after the last (and currently only) sheet so it
# lands at the end of the tab strip, then give it its final name.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "productimages"

# --- Column widths (characters) ---
$ws2.Columns.Item(1).ColumnWidth = 13.16
$ws2.Columns.Item(2).ColumnWidth = 33.65
$ws2.Columns.Item(3).ColumnWidth = 13.0

# --- Formatting: clone the "register" sheet's header style (border + wrap +
#     yellow fill), then bump it to bold for this sheet's header row. ---
$ws1.Range("A1:C1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)
$ws2.Range("A1:C1").Font.Bold = $true
$ws2.Rows.Item(1).RowHeight = 30

# --- Formatting: clone the plain bordered/wrapped body style for the data
#     rows (A2:B4 stay regular text; C2:C4 get the quote-prefixed style). ---
$ws1.Range("A2").Copy()
$ws2.Range("A2:B4").PasteSpecial(-4122)
$ws2.Range("C2:C4").PasteSpecial(-4122)

# --- Header values ---
$ws2.Range("A1").Value = "searchitem"
$ws2.Range("B1").Value = "products"
$ws2.Range("C1").Value = "imagecount"

# --- Data values (write order intentionally matches the original author's
#     entry order so the shared-string table indices line up). ---
$ws2.Range("B2").Value = "Samsung Galaxy Tab 10.1"
$ws2.Range("A2").Value = "Samsung"
$ws2.Range("A3").Value = "Macbook"
$ws2.Range("B3").Value = "MacBook Air"
$ws2.Range("A4").Value = "iMac"
$ws2.Range("B4").Value = "iMac"

# Image counts are entered with a leading apostrophe so they are stored as
# text (quote-prefixed), matching the source data.
$ws2.Range("C2").Value = "'7"
$ws2.Range("C3").Value = "'4"
$ws2.Range("C4").Value = "'3"

# Make the new sheet the active tab with C5 selected (next empty row below
# the data), mirroring where the author's cursor was left.
[void]$ws2.Activate()
[void]$ws2.Range("C5").Select()
